$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-4 are re-sorted by date (column D) ascending.
# Row 5 already holds the earliest date and is left untouched.
# New row 2 = old row 4, new row 3 = old row 2, new row 4 = old row 3.

$ws.Range("D2").Value = 44280
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 500

$ws.Range("D3").Value = 44284
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 500

$ws.Range("D4").Value = 44291
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550
